# ---------------------------------------------------------------------------
# Project 2 phase 2 update:
#  - Snapshot the current "Data" sheet into a new "Old Data" sheet (so the
#    old/pre-update numbers are preserved for comparison), inserted between
#    "Data" and "Precision-Recall Graph".
#  - Add a difference table (rows 15-24) on "Old Data" comparing Data vs
#    Old Data, with conditional formatting (green = improved, red = worse).
#  - Update the "Data" sheet numbers themselves to the new phase-2 results.
#  - Freeze the header rows on "Data" and move the selection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

# 1) Duplicate the "Data" sheet (carries over all formatting/styles/values)
#    and place it right after "Data" -- this becomes "Old Data" and keeps
#    the pre-edit numbers.
$dataSheet.Copy($null, $dataSheet)
$oldData = $wb.Worksheets.Item(2)
$oldData.Name = "Old Data"

# A locally scoped defined name "ir" pointing at the Old Data range, mirroring
# the one that already exists on "Data".
$oldData.Names.Add("ir", "='Old Data'!`$A`$3:`$B`$13")

# Reset Old Data's view (no frozen panes there, just a plain selection).
$oldData.Application.ActiveWindow.FreezePanes = $false
$oldData.Range("E3").Select()

# 2) Add the diff rows (15-24) under the copied table on "Old Data".
$diffRows = 3..12
$destRow = 15
foreach ($srcRow in $diffRows) {
    foreach ($col in @("B","C","D","E","F","G","H","I")) {
        $cell = $oldData.Range($col + $destRow)
        $cell.Formula = "=Data!" + $col + $srcRow + "-'Old Data'!" + $col + $srcRow
    }
    $destRow = $destRow + 1
}

# Style the diff cells with the numeric style used elsewhere (index 3 in the
# original styles.xml: 0.000000000000000 format, centered).
$oldData.Range("B15:I24").NumberFormat = "0.000000000000000"
$oldData.Range("B15:I24").HorizontalAlignment = -4108

# Conditional formatting: green when improved (>0), red when worse (<0).
$ws_range = $oldData.Range("B15:I24")
$ws_range.FormatConditions.Delete()

$fc1 = $ws_range.FormatConditions.Add(1, 5, "0")
$fc1.Font.Color = 3877552
$fc1.Interior.Color = 13561798

$fc2 = $ws_range.FormatConditions.Add(1, 6, "0")
$fc2.Font.Color = 402
$fc2.Interior.Color = 13551615

# 3) Update the "Data" sheet's own numbers to the new phase-2 results.
$dataSheet.Range("E3").Value = 0.73391980507327703
$dataSheet.Range("I3").Value = 0.67689392373980894

$dataSheet.Range("E4").Value = 0.73191980507327703
$dataSheet.Range("I4").Value = 0.675075741921627

$dataSheet.Range("E5").Value = 0.72358647173994395
$dataSheet.Range("I5").Value = 0.667075741921628

$dataSheet.Range("D6").Value = 0.70097670182746696
$dataSheet.Range("E6").Value = 0.68830869396216598
$dataSheet.Range("I6").Value = 0.64585351969940497

$dataSheet.Range("D7").Value = 0.69364336849413299
$dataSheet.Range("E7").Value = 0.68230869396216598
$dataSheet.Range("I7").Value = 0.620569969915855

$dataSheet.Range("D8").Value = 0.68896490685915002
$dataSheet.Range("E8").Value = 0.67626023949265901
$dataSheet.Range("I8").Value = 0.60983187467776001

$dataSheet.Range("D9").Value = 0.58237497946083505
$dataSheet.Range("E9").Value = 0.56727331416473703
$dataSheet.Range("I9").Value = 0.53920118172723897

$dataSheet.Range("D10").Value = 0.54746383413609201
$dataSheet.Range("E10").Value = 0.53469613578090702
$dataSheet.Range("I10").Value = 0.51887298378139401

$dataSheet.Range("D11").Value = 0.53059166107807798
$dataSheet.Range("E11").Value = 0.51792053665236704
$dataSheet.Range("I11").Value = 0.490916187223364

$dataSheet.Range("D12").Value = 0.50564303483807505
$dataSheet.Range("E12").Value = 0.49169667788577898
$dataSheet.Range("I12").Value = 0.46636176594147899

$dataSheet.Range("D13").Value = 0.50546781581429501
$dataSheet.Range("E13").Value = 0.49169667788577898
$dataSheet.Range("I13").Value = 0.465121999859608

# 4) Freeze header rows on "Data" and move the active selection/pane.
$dataSheet.Activate()
$dataSheet.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$dataSheet.Range("C3").Select()

$wb.Save()
